$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Commit: "added one json for time bucket analysis"
#
# The underlying data-generation pipeline re-ran after a new JSON source was
# folded into the time-bucket analysis, which re-ordered the material rows
# of the table (row 2 <-> row 3 swap their whole record, rows 4-6 rotate:
# 5->4, 6->5, 4->6) while the header row (row 1) and the overall shape
# (A1:E6) stay the same. Capture the final values/hyperlinks for every data
# row explicitly so every cell lands on its correct final content.
# ---------------------------------------------------------------------------

# New row 2  (was old row 3 - "His Life and Times ...")
$ws.Range("A2").Value = "His Life and Times: Albert E. Castel: 9780806130811: Amazon.com: Books"
$ws.Range("B2").Value = "2015-07-22T00:00:00UTC"
$ws.Range("C2").Value = 55487
$ws.Range("D2").Value = "day_31_beyond"
$ws.Range("E2").Value = "https://www.amazon.com/William-Clarke-Quantrill-Life-Times/dp/0806130814/"

# New row 3 (was old row 2 - "63 Roster ...")
$ws.Range("A3").Value = "63 Roster of of the Victims of the Lawrence Raid Published on the 147th Anniversary."
$ws.Range("B3").Value = "2020-09-27T00:00:00UTC"
$ws.Range("C3").Value = 57381
$ws.Range("D3").Value = "day_31_beyond"
$ws.Range("E3").Value = "http://quantrillsguerrillas.com/en/articles/126-t-charles-edwin-wells.html"

# New row 4 (was old row 5 - "Guerilla Operations ...")
$ws.Range("A4").Value = "Guerilla Operations in the Trans-Mississippi"
$ws.Range("B4").Value = "1-01-01T00:00:00UTC"
$ws.Range("C4").Value = "unknown"
$ws.Range("D4").Value = "unknown"
$ws.Range("E4").Value = "http://stellar-one.com/guerilla/"

# New row 5 (was old row 6 - "H.M. Simpson ...")
$ws.Range("A5").Value = "H.M. Simpson to Hiram Hill"
$ws.Range("B5").Value = "1-01-01T00:00:00UTC"
$ws.Range("C5").Value = "unknown"
$ws.Range("D5").Value = "unknown"
$ws.Range("E5").Value = "http://www.kansasmemory.org/item/213271"

# New row 6 (was old row 4 - "Quantrill and the border wars")
$ws.Range("A6").Value = "Quantrill and the border wars"
$ws.Range("B6").Value = "1-01-01T00:00:00UTC"
$ws.Range("C6").Value = "unknown"
$ws.Range("D6").Value = "unknown"
$ws.Range("E6").Value = "https://archive.org/details/quantrillborderw00connuoft"

# Rebuild the "uri" column hyperlinks so each moved row keeps its link
# (the link target always mirrors the displayed URI text in column E).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), $ws.Range("E2").Value2)
$ws.Hyperlinks.Add($ws.Range("E3"), $ws.Range("E3").Value2)
$ws.Hyperlinks.Add($ws.Range("E4"), $ws.Range("E4").Value2)
$ws.Hyperlinks.Add($ws.Range("E5"), $ws.Range("E5").Value2)
$ws.Hyperlinks.Add($ws.Range("E6"), $ws.Range("E6").Value2)

# Hyperlinks.Add() mints a brand-new (functionally identical) cell style
# instead of reusing the workbook's existing "Hyperlink" cell style - put
# E2:E6 back on the original named style so formatting is unchanged.
$ws.Range("E2:E6").Style = "Hyperlink"
